$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.993.49"
Set-TextValue "E2" "  -0.42%  "

Set-TextValue "D3" "1.880.11"
Set-TextValue "E3" "  -1.59%  "

Set-TextValue "D4" "0.9993"
Set-TextValue "E4" "  -0.15%  "

Set-TextValue "D5" "243.03"
Set-TextValue "E5" "  -3.36%  "

Set-TextValue "D6" "0.9992"
Set-TextValue "E6" "  -0.19%  "

Set-TextValue "D7" "0.4966"
Set-TextValue "E7" "  -2.65%  "

Set-TextValue "D8" "0.2923"
Set-TextValue "E8" "  -1.05%  "

Set-TextValue "D9" "0.06650"
Set-TextValue "E9" "  -1.83%  "

Set-TextValue "D10" "1.879.58"
Set-TextValue "E10" "  -1.58%  "

Set-TextValue "D11" "16.81"
Set-TextValue "E11" "  -2.58%  "

Set-TextValue "D12" "0.07246"
Set-TextValue "E12" "  -1.62%  "

Set-TextValue "D13" "0.6682"
Set-TextValue "E13" "  -3.05%  "

Set-TextValue "D14" "86.62"
Set-TextValue "E14" "  +0.16%  "

Set-TextValue "D15" "4.926"
Set-TextValue "E15" "  +1.08%  "

Set-TextValue "D16" "29.974.73"
Set-TextValue "E16" "  -0.51%  "

Set-TextValue "D17" "0.000007876"
Set-TextValue "E17" "  -3.05%  "

Set-TextValue "D18" "0.9986"
Set-TextValue "E18" "  -0.24%  "

Set-TextValue "D19" "12.81"
Set-TextValue "E19" "  -1.23%  "

Set-TextValue "D20" "2.123.62"
Set-TextValue "E20" "  -1.62%  "

Set-TextValue "D21" "0.9991"
Set-TextValue "E21" "  -0.12%  "

Set-TextValue "D22" "4.780"
Set-TextValue "E22" "  -0.98%  "

Set-TextValue "D23" "5.762"
Set-TextValue "E23" "  +0.46%  "

Set-TextValue "D24" "9.087"
Set-TextValue "E24" "  -0.59%  "

Set-TextValue "D25" "142.82"
Set-TextValue "E25" "  +5.82%  "

Set-TextValue "D26" "149.70"
Set-TextValue "E26" "  +1.91%  "

Set-TextValue "D27" "17.06"
Set-TextValue "E27" "  -0.11%  "

Set-TextValue "D28" "1.918"
Set-TextValue "E28" "  -3.63%  "

Set-TextValue "E29" "  +0.05%  "

Set-TextValue "D30" "4.213"
Set-TextValue "E30" "  -0.31%  "

Set-TextValue "D31" "0.08770"
Set-TextValue "E31" "  -0.15%  "

Set-TextValue "D32" "3.971"
Set-TextValue "E32" "  -0.92%  "

Set-TextValue "D33" "0.05093"
Set-TextValue "E33" "  +0.67%  "

Set-TextValue "B34" "ARBITRUM"
Set-TextValue "C34" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D34" "1.120"
Set-TextValue "E34" "  -2.06%  "

Set-TextValue "B35" "ImmutableX"
Set-TextValue "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D35" "0.7140"
Set-TextValue "E35" "  +0.09%  "

Set-TextValue "D36" "2.668"

Set-TextValue "D37" "0.01805"
Set-TextValue "E37" "  +6.66%  "

Set-TextValue "D38" "2.692"
Set-TextValue "E38" "  -4.06%  "

Set-TextValue "D39" "2.176"
Set-TextValue "E39" "  -4.48%  "

Set-TextValue "D40" "0.9347"
Set-TextValue "E40" "  -3.52%  "

Set-TextValue "B41" "FraxShare"
Set-TextValue "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "5.809"
Set-TextValue "E41" "  -5.29%  "

Set-TextValue "B42" "TheSandbox"
Set-TextValue "C42" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D42" "0.4266"
Set-TextValue "E42" "  -0.50%  "

Set-TextValue "D43" "0.9986"
Set-TextValue "E43" "  -0.03%  "

Set-TextValue "D44" "102.68"
Set-TextValue "E44" "  -1.99%  "

Set-TextValue "D45" "7.456"

Set-TextValue "D46" "0.1271"
Set-TextValue "E46" "  -0.49%  "

Set-TextValue "E47" "  -1.27%  "

Set-TextValue "D48" "32.64"
Set-TextValue "E48" "  -1.37%  "

Set-TextValue "D49" "0.3796"
Set-TextValue "E49" "  +0.01%  "

Set-TextValue "D50" "8.307"
Set-TextValue "E50" "  -1.48%  "

Set-TextValue "D51" "56.04"
Set-TextValue "E51" "  -1.27%  "
